# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.907.40'
$ws.Range("E2").Value = '  +8.10%  '
$ws.Range("D3").Value = '1.813.34'
$ws.Range("E3").Value = '  +4.97%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = "'246.44"
$ws.Range("E5").Value = '  +2.45%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D7").Value = "'0.4927"
$ws.Range("E7").Value = '  +1.94%  '
$ws.Range("D8").Value = "'43.89"
$ws.Range("E8").Value = '  +6.54%  '
$ws.Range("D9").Value = "'0.2778"
$ws.Range("E9").Value = '  +7.39%  '
$ws.Range("D10").Value = "'0.06399"
$ws.Range("E10").Value = '  +3.44%  '
$ws.Range("D11").Value = '1.810.33'
$ws.Range("E11").Value = '  +4.82%  '
$ws.Range("D12").Value = "'16.73"
$ws.Range("E12").Value = '  +5.34%  '
$ws.Range("D13").Value = "'0.07075"
$ws.Range("E13").Value = '  +2.99%  '
$ws.Range("D14").Value = "'0.6435"
$ws.Range("E14").Value = '  +6.45%  '
$ws.Range("D15").Value = "'83.91"
$ws.Range("E15").Value = '  +8.98%  '
$ws.Range("D16").Value = "'4.675"
$ws.Range("E16").Value = '  +4.57%  '
$ws.Range("D17").Value = '28.932.96'
$ws.Range("E17").Value = '  +8.93%  '
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = '  +0.12%  '
$ws.Range("D19").Value = "'0.000007313"
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = "'12.24"
$ws.Range("E21").Value = '  +7.68%  '
$ws.Range("D22").Value = '2.042.46'
$ws.Range("E22").Value = '  +4.81%  '
$ws.Range("D23").Value = "'4.572"
$ws.Range("E23").Value = '  +3.43%  '
$ws.Range("D24").Value = "'8.799"
$ws.Range("E24").Value = '  +2.94%  '
$ws.Range("D25").Value = "'5.339"
$ws.Range("E25").Value = '  +5.57%  '
$ws.Range("D26").Value = "'143.33"
$ws.Range("E26").Value = '  +2.46%  '
$ws.Range("D27").Value = "'128.86"
$ws.Range("E27").Value = '  +20.94%  '
$ws.Range("D28").Value = "'16.46"
$ws.Range("E28").Value = '  +8.09%  '
$ws.Range("D29").Value = "'1.883"
$ws.Range("E29").Value = '  +5.91%  '
$ws.Range("D30").Value = "'1.411"
$ws.Range("E30").Value = '  +3.17%  '
$ws.Range("D31").Value = "'4.128"
$ws.Range("E31").Value = '  +2.90%  '
$ws.Range("D32").Value = "'0.08340"
$ws.Range("E32").Value = '  +5.07%  '
$ws.Range("D33").Value = "'3.777"
$ws.Range("E33").Value = '  +2.90%  '
$ws.Range("D34").Value = "'0.04940"
$ws.Range("E34").Value = '  +9.41%  '
$ws.Range("E35").Value = '  +9.56%  '
$ws.Range("D36").Value = "'2.698"
$ws.Range("E36").Value = '  +3.86%  '
$ws.Range("D37").Value = "'0.6714"
$ws.Range("E37").Value = '  +8.60%  '
$ws.Range("E38").Value = '  +14.40%  '
$ws.Range("D39").Value = "'2.744"
$ws.Range("E39").Value = '  +11.71%  '
$ws.Range("D40").Value = "'0.9537"
$ws.Range("E40").Value = '  +1.98%  '
$ws.Range("D41").Value = "'6.130"
$ws.Range("E41").Value = '  +9.24%  '
$ws.Range("D42").Value = "'0.01585"
$ws.Range("E42").Value = '  +5.72%  '
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("D44").Value = "'100.78"
$ws.Range("E44").Value = '  +0.99%  '
$ws.Range("D45").Value = "'0.4079"
$ws.Range("E45").Value = '  +6.42%  '
$ws.Range("D46").Value = "'7.125"
$ws.Range("E46").Value = '  +4.95%  '
$ws.Range("D47").Value = "'0.1220"
$ws.Range("E47").Value = '  +5.58%  '
$ws.Range("D48").Value = "'0.05525"
$ws.Range("E48").Value = '  +3.10%  '
$ws.Range("D49").Value = "'8.138"
$ws.Range("E49").Value = '  +3.02%  '
$ws.Range("E50").Value = '  +5.00%  '
$ws.Range("D51").Value = "'0.3614"
$ws.Range("E51").Value = '  +7.75%  '
